$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Insert a new row at 64 (shifts rows 64.. down by one), matching the
# addition of Santoker Q/R Series support documented in the commit message.
$ws.Rows.Item(64).Insert()

$ws.Range("B64").Value = "santoker(<target>,<value>)"
$ws.Range("C64").Value = "sends integer <value> to <target> register specified by as byte in hex notation like “fa” via the Santoker Network protocol"

# Update selection/active cell to reflect the new row position.
$null = $ws.Range("B64:C64").Select()
